# إضافة حدث جديد في Card22
# Adds a new service-log row (row 15) to the "Card22" worksheet and fills
# the previously-blank placeholder cells (B14:K14) on the prior row with
# the literal text "nan", matching the source data export's convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card22")

# --- Fill the (until now empty) placeholder cells on row 14 with "nan" ---
$placeholderCols = @("B","C","D","E","F","G","H","I","J","K")
foreach ($col in $placeholderCols) {
    $ws.Range($col + "14").Value = "nan"
}

# --- Append the new event as row 15 ---
# Column A holds the card number ("22") as text, matching the rest of the
# column. A plain Range.Value assignment of a numeric-looking string gets
# auto-coerced to a number by Excel, so we compute it as text via TEXT()
# and paste-special the value back in, which keeps it a text cell without
# touching the cell's style.
$ws.Range("Z1").Formula = "=TEXT(22,""0"")"
$ws.Range("Z1").Copy()
$ws.Range("A15").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

# B15:K15 stay blank, same as the previous template row used to be.

$ws.Range("L15").Value = "14\10\2024"
# M15 stays blank.
$ws.Range("N15").Value = "تم تغيير السستم من ax اليay"
$ws.Range("O15").Value = "تيم العمل"
